$p = $ppt.ActivePresentation

# Slide 4: Title "Option 2: Using 24-Bit Session ID Field" -> "Option 2: 24-bit Session ID Field"
$s4 = $p.Slides.Item(4)
$title4 = $s4.Shapes.Item(1).TextFrame.TextRange
$title4.Paragraphs(1, 1).Runs(1, 1).Text = "Option 2: 24-bit Session ID Field"

# Slide 6: Title "Option 4: 16-bit Session ID - 1" -> "Option 4: 16-bit Session ID Field - 1"
$s6 = $p.Slides.Item(6)
$title6 = $s6.Shapes.Item(1).TextFrame.TextRange
$title6.Paragraphs(1, 1).Runs(1, 1).Text = "Option 4: 16-bit Session ID Field - 1"

# Slide 6: TextBox first bullet "Session ID has only 16-bits - is this enough?"
#          -> "Session ID has only 16-bits - is this good enough?"
$tb6 = $s6.Shapes.Item(6)
$tb6.TextFrame.TextRange.Paragraphs(1, 1).Runs(1, 1).Text = "Session ID has only 16-bits – is this good enough?"

# Slide 7: Title "Option 5: 16-bit Session ID - 2" -> "Option 5: 16-bit Session ID Field - 2"
$s7 = $p.Slides.Item(7)
$title7 = $s7.Shapes.Item(1).TextFrame.TextRange
$title7.Paragraphs(1, 1).Runs(1, 1).Text = "Option 5: 16-bit Session ID Field - 2"

# Slide 7: TextBox first bullet "Session ID has only 16-bits - is this enough?"
#          -> "Session ID has only 16-bits - is this good enough?"
$tb7 = $s7.Shapes.Item(6)
$tb7.TextFrame.TextRange.Paragraphs(1, 1).Runs(1, 1).Text = "Session ID has only 16-bits – is this good enough?"
